$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2550
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3968

$ws.Range("H76").Value = 2500
$ws.Range("I76").Value = 2000
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 2000
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -1685
$ws.Range("N76").Value = -3630

$ws.Range("H79").Value = 2500
$ws.Range("I79").Value = 2000
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 2000
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -908
$ws.Range("N79").Value = -5184

$ws.Range("H93").Value = 65300.5
$ws.Range("J93").Value = 65300.5
$ws.Range("L93").Value = 65300.5
$ws.Range("N93").Value = -70292.5

$ws.Range("H98").Value = 626.0606
$ws.Range("I98").Value = 657.2414
$ws.Range("K98").Value = 657.2414
$ws.Range("M98").Value = 840.7586

$ws.Range("H113").Value = 48269.84
$ws.Range("I113").Value = 75303.36
$ws.Range("J113").Value = 13863.546
$ws.Range("K113").Value = 75303.36
$ws.Range("L113").Value = 13863.546
$ws.Range("M113").Value = -72049.36
$ws.Range("N113").Value = -20371.546

$ws.Range("H122").Value = 626.0606
$ws.Range("I122").Value = 657.2414
$ws.Range("K122").Value = 1971.7242
$ws.Range("M122").Value = 478.2757999999999

$ws.Range("H132").Value = 23953890
$ws.Range("I132").Value = 25719574
$ws.Range("J132").Value = 1000000
$ws.Range("K132").Value = 77158722
$ws.Range("L132").Value = 3000000
$ws.Range("M132").Value = -77156192
$ws.Range("N132").Value = -3005060

$ws.Range("H135").Value = 743.38464
$ws.Range("I135").Value = 472
$ws.Range("K135").Value = 4248
$ws.Range("M135").Value = -1713

$ws.Range("H138").Value = 3995.2896
$ws.Range("I138").Value = 2481.7
$ws.Range("J138").Value = 4224.621
$ws.Range("K138").Value = 7445.099999999999
$ws.Range("L138").Value = 12673.863
$ws.Range("M138").Value = -2305.099999999999
$ws.Range("N138").Value = -22953.863

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3379.3125
$ws.Range("I2").Value = 3204.9
$ws.Range("J2").Value = 3670
$ws.Range("K2").Value = 3204.9
$ws.Range("L2").Value = 3670
$ws.Range("M2").Value = -3091.9
$ws.Range("N2").Value = -3896

$ws.Range("H32").Value = 2662.1177
$ws.Range("I32").Value = 2339.9583
$ws.Range("K32").Value = 2339.9583
$ws.Range("M32").Value = -2052.9583

$ws.Range("H116").Value = 3379.3125
$ws.Range("I116").Value = 3204.9
$ws.Range("J116").Value = 3670
$ws.Range("K116").Value = 3204.9
$ws.Range("L116").Value = 3670
$ws.Range("M116").Value = -910.9000000000001
$ws.Range("N116").Value = -8258

$ws.Range("H122").Value = 2067.1875
$ws.Range("I122").Value = 1713.7273
$ws.Range("J122").Value = 2844.8
$ws.Range("K122").Value = 5141.1819
$ws.Range("L122").Value = 8534.400000000001
$ws.Range("M122").Value = -2691.1819
$ws.Range("N122").Value = -13434.4

$ws.Range("H132").Value = 2795
$ws.Range("I132").Value = 2540.647
$ws.Range("K132").Value = 7621.941
$ws.Range("M132").Value = -5091.941

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3379.3125
$ws.Range("I3").Value = 3204.9
$ws.Range("J3").Value = 3670
$ws.Range("K3").Value = 3204.9
$ws.Range("L3").Value = 3670
$ws.Range("M3").Value = -3090.9
$ws.Range("N3").Value = -3898

$ws.Range("H94").Value = 1817.3043
$ws.Range("I94").Value = 1637.8
$ws.Range("K94").Value = 1637.8
$ws.Range("M94").Value = -1186.8

$ws.Range("H99").Value = 4269.9565
$ws.Range("I99").Value = 1318.2941
$ws.Range("J99").Value = 12633
$ws.Range("K99").Value = 1318.2941
$ws.Range("L99").Value = 12633
$ws.Range("M99").Value = 179.7058999999999
$ws.Range("N99").Value = -15629

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27482.463
$ws.Range("I31").Value = 28728.838
$ws.Range("K31").Value = 28728.838
$ws.Range("M31").Value = -28433.838

$ws.Range("H34").Value = 27482.463
$ws.Range("I34").Value = 28728.838
$ws.Range("K34").Value = 28728.838
$ws.Range("M34").Value = -28526.838

$ws.Range("H58").Value = 1813.4736
$ws.Range("I58").Value = 1840.8823
$ws.Range("K58").Value = 1840.8823
$ws.Range("M58").Value = -1637.8823

$ws.Range("H98").Value = 200000
$ws.Range("J98").Value = 200000
$ws.Range("L98").Value = 200000
$ws.Range("N98").Value = -204492

$ws.Range("H107").Value = 1026.25
$ws.Range("J107").Value = 623.5
$ws.Range("L107").Value = 623.5
$ws.Range("N107").Value = -4463.5

$ws.Range("H134").Value = 25340
$ws.Range("I134").Value = 10542.941
$ws.Range("K134").Value = 31628.823
$ws.Range("M134").Value = -29093.823

$ws.Range("H136").Value = 1813.4736
$ws.Range("I136").Value = 1840.8823
$ws.Range("K136").Value = 5522.6469
$ws.Range("M136").Value = -2972.6469

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 5835.3335
$ws.Range("I116").Value = 2503.6
$ws.Range("K116").Value = 7510.799999999999
$ws.Range("M116").Value = -4068.799999999999

$ws.Range("H121").Value = 6363.4443
$ws.Range("I121").Value = 4261.1875
$ws.Range("J121").Value = 8045.25
$ws.Range("K121").Value = 12783.5625
$ws.Range("L121").Value = 24135.75
$ws.Range("M121").Value = -11473.5625
$ws.Range("N121").Value = -26755.75

$ws.Range("H128").Value = 349992.72
$ws.Range("I128").Value = 349992.72
$ws.Range("K128").Value = 1049978.16
$ws.Range("M128").Value = -1044998.16

$ws.Range("H138").Value = 12858.333
$ws.Range("I138").Value = 11404.167
$ws.Range("K138").Value = 34212.501
$ws.Range("M138").Value = -29072.501

$ws.Range("H140").Value = 11066
$ws.Range("I140").Value = 11066
$ws.Range("K140").Value = 33198
$ws.Range("M140").Value = -28018

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9652.210999999999
$ws.Range("I70").Value = 9317.177
$ws.Range("J70").Value = 12500
$ws.Range("K70").Value = 9317.177
$ws.Range("L70").Value = 12500
$ws.Range("M70").Value = -9047.177
$ws.Range("N70").Value = -13040

$ws.Range("H73").Value = 9652.210999999999
$ws.Range("I73").Value = 9317.177
$ws.Range("J73").Value = 12500
$ws.Range("K73").Value = 9317.177
$ws.Range("L73").Value = 12500
$ws.Range("M73").Value = -8381.177
$ws.Range("N73").Value = -14372

$ws.Range("H102").Value = 1592.9286
$ws.Range("I102").Value = 1578.6578
$ws.Range("K102").Value = 1578.6578
$ws.Range("M102").Value = 43.34220000000005

$ws.Range("H122").Value = 1880.4
$ws.Range("I122").Value = 1499.75
$ws.Range("J122").Value = 2134.1667
$ws.Range("K122").Value = 4499.25
$ws.Range("L122").Value = 6402.500100000001
$ws.Range("M122").Value = -2049.25
$ws.Range("N122").Value = -11302.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2164.5881
$ws.Range("I22").Value = 1868
$ws.Range("J22").Value = 2326.3635
$ws.Range("K22").Value = 1868
$ws.Range("L22").Value = 2326.3635
$ws.Range("M22").Value = -1573
$ws.Range("N22").Value = -2916.3635

$ws.Range("H27").Value = 2164.5881
$ws.Range("I27").Value = 1868
$ws.Range("J27").Value = 2326.3635
$ws.Range("K27").Value = 1868
$ws.Range("L27").Value = 2326.3635
$ws.Range("M27").Value = -1761
$ws.Range("N27").Value = -2540.3635

$ws.Range("H98").Value = 41000
$ws.Range("J98").Value = 41000
$ws.Range("L98").Value = 41000
$ws.Range("N98").Value = -46990

$ws.Range("H122").Value = 5825.577
$ws.Range("I122").Value = 5480.75
$ws.Range("J122").Value = 6975
$ws.Range("K122").Value = 16442.25
$ws.Range("L122").Value = 20925
$ws.Range("M122").Value = -13992.25
$ws.Range("N122").Value = -25825

$ws.Range("H132").Value = 5152.05
$ws.Range("I132").Value = 4922.3076
$ws.Range("K132").Value = 14766.9228
$ws.Range("M132").Value = -12236.9228

$ws.Range("H136").Value = 6350.9
$ws.Range("I136").Value = 6722.6665
$ws.Range("J136").Value = 3005
$ws.Range("K136").Value = 20167.9995
$ws.Range("L136").Value = 9015
$ws.Range("M136").Value = -17617.9995
$ws.Range("N136").Value = -14115

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2777.5
$ws.Range("I122").Value = 2455.25
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 7365.75
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -4915.75
$ws.Range("N122").Value = -22900

$ws.Range("H126").Value = 3638.6155
$ws.Range("I126").Value = 3557.0625
$ws.Range("J126").Value = 4011.4285
$ws.Range("K126").Value = 10671.1875
$ws.Range("L126").Value = 12034.2855
$ws.Range("M126").Value = -8201.1875
$ws.Range("N126").Value = -16974.2855

$ws.Range("H132").Value = 2181.516
$ws.Range("I132").Value = 1986.7858
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 5960.357400000001
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -3430.357400000001
$ws.Range("N132").Value = -17057

$ws.Range("H136").Value = 2675.7856
$ws.Range("I136").Value = 2607.6667
$ws.Range("K136").Value = 7823.000100000001
$ws.Range("M136").Value = -5273.000100000001
